# Scheduled-runner style refresh of market/profit figures in the
# Chocobo_Profits workbook. Each worksheet (one per crafting class) holds a
# leve table whose H:N columns are live price/profit data pulled from an
# external market-board API. This script overwrites those columns with the
# newly refreshed values for the handful of rows whose source data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1225.75
$ws.Range("I28").Value = 1265.5883
$ws.Range("K28").Value = 1265.5883
$ws.Range("M28").Value = -780.5882999999999
$ws.Range("H43").Value = 2040.6316
$ws.Range("I43").Value = 956.6667
$ws.Range("J43").Value = 3898.8572
$ws.Range("K43").Value = 956.6667
$ws.Range("L43").Value = 3898.8572
$ws.Range("M43").Value = -887.6667
$ws.Range("N43").Value = -4036.8572
$ws.Range("H62").Value = 1281.5
$ws.Range("I62").Value = 1137.8
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 1137.8
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -513.8
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 1281.5
$ws.Range("I65").Value = 1137.8
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 5689
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -2569
$ws.Range("N65").Value = -16240
$ws.Range("H132").Value = 36216270
$ws.Range("I132").Value = 55775150
$ws.Range("J132").Value = 1010280.2
$ws.Range("K132").Value = 167325450
$ws.Range("L132").Value = 3030840.6
$ws.Range("M132").Value = -167322920
$ws.Range("N132").Value = -3035900.6
$ws.Range("H135").Value = 1163.5264
$ws.Range("I135").Value = 693.13336
$ws.Range("J135").Value = 2927.5
$ws.Range("K135").Value = 6238.20024
$ws.Range("L135").Value = 26347.5
$ws.Range("M135").Value = -3703.20024
$ws.Range("N135").Value = -31417.5
$ws.Range("H136").Value = 49118.168
$ws.Range("J136").Value = 49118.168
$ws.Range("L136").Value = 49118.168
$ws.Range("N136").Value = -59318.168
$ws.Range("H137").Value = 1038993.56
$ws.Range("I137").Value = 3179654
$ws.Range("J137").Value = 3190.1292
$ws.Range("K137").Value = 9538962
$ws.Range("L137").Value = 9570.3876
$ws.Range("M137").Value = -9536412
$ws.Range("N137").Value = -14670.3876

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1473.625
$ws.Range("I2").Value = 1486.96
$ws.Range("K2").Value = 1486.96
$ws.Range("M2").Value = -1373.96
$ws.Range("H45").Value = 3281.5833
$ws.Range("I45").Value = 3100
$ws.Range("J45").Value = 3372.375
$ws.Range("K45").Value = 3100
$ws.Range("L45").Value = 3372.375
$ws.Range("M45").Value = -2723
$ws.Range("N45").Value = -4126.375
$ws.Range("H74").Value = 486220.62
$ws.Range("I74").Value = 820733.0600000001
$ws.Range("J74").Value = 3036
$ws.Range("K74").Value = 820733.0600000001
$ws.Range("L74").Value = 3036
$ws.Range("M74").Value = -819859.0600000001
$ws.Range("N74").Value = -4784
$ws.Range("H77").Value = 486220.62
$ws.Range("I77").Value = 820733.0600000001
$ws.Range("J77").Value = 3036
$ws.Range("K77").Value = 4103665.3
$ws.Range("L77").Value = 15180
$ws.Range("M77").Value = -4099297.3
$ws.Range("N77").Value = -23916
$ws.Range("H116").Value = 1473.625
$ws.Range("I116").Value = 1486.96
$ws.Range("K116").Value = 1486.96
$ws.Range("M116").Value = 807.04
$ws.Range("H122").Value = 2150.8462
$ws.Range("I122").Value = 1286.1
$ws.Range("K122").Value = 3858.3
$ws.Range("M122").Value = -1408.3
$ws.Range("H132").Value = 2441.7778
$ws.Range("I132").Value = 1484.2858
$ws.Range("J132").Value = 4018.8235
$ws.Range("K132").Value = 4452.857400000001
$ws.Range("L132").Value = 12056.4705
$ws.Range("M132").Value = -1922.857400000001
$ws.Range("N132").Value = -17116.4705

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1963.125
$ws.Range("I99").Value = 1770
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 1770
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = -272
$ws.Range("N99").Value = -5796
$ws.Range("H107").Value = 2760.4
$ws.Range("I107").Value = 3934
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 3934
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -2014
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2429.196
$ws.Range("I31").Value = 946.625
$ws.Range("J31").Value = 3747.037
$ws.Range("K31").Value = 946.625
$ws.Range("L31").Value = 3747.037
$ws.Range("M31").Value = -651.625
$ws.Range("N31").Value = -4337.037
$ws.Range("H34").Value = 2429.196
$ws.Range("I34").Value = 946.625
$ws.Range("J34").Value = 3747.037
$ws.Range("K34").Value = 946.625
$ws.Range("L34").Value = 3747.037
$ws.Range("M34").Value = -744.625
$ws.Range("N34").Value = -4151.037
$ws.Range("H58").Value = 2626.5
$ws.Range("J58").Value = 5068.9
$ws.Range("L58").Value = 5068.9
$ws.Range("N58").Value = -5474.9
$ws.Range("H122").Value = 6000
$ws.Range("J122").Value = 10000
$ws.Range("L122").Value = 30000
$ws.Range("N122").Value = -34900
$ws.Range("H132").Value = 3973.353
$ws.Range("I132").Value = 3284.1875
$ws.Range("K132").Value = 9852.5625
$ws.Range("M132").Value = -7322.5625
$ws.Range("H134").Value = 2664.9
$ws.Range("I134").Value = 1206.125
$ws.Range("J134").Value = 8500
$ws.Range("K134").Value = 3618.375
$ws.Range("L134").Value = 25500
$ws.Range("M134").Value = -1083.375
$ws.Range("N134").Value = -30570
$ws.Range("H136").Value = 2626.5
$ws.Range("J136").Value = 5068.9
$ws.Range("L136").Value = 15206.7
$ws.Range("N136").Value = -20306.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 770.9400000000001
$ws.Range("J131").Value = 800.10986
$ws.Range("L131").Value = 2400.32958
$ws.Range("N131").Value = -12480.32958

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2770.2
$ws.Range("I7").Value = 2037.5
$ws.Range("J7").Value = 5701
$ws.Range("K7").Value = 2037.5
$ws.Range("L7").Value = 5701
$ws.Range("M7").Value = -1925.5
$ws.Range("N7").Value = -5925
$ws.Range("H40").Value = 5749.8887
$ws.Range("I40").Value = 4968.6924
$ws.Range("K40").Value = 4968.6924
$ws.Range("M40").Value = -4832.6924
$ws.Range("H126").Value = 2770.2
$ws.Range("I126").Value = 2037.5
$ws.Range("J126").Value = 5701
$ws.Range("K126").Value = 6112.5
$ws.Range("L126").Value = 17103
$ws.Range("M126").Value = -3642.5
$ws.Range("N126").Value = -22043
$ws.Range("H132").Value = 7546.1816
$ws.Range("I132").Value = 7402.6665
$ws.Range("K132").Value = 22207.9995
$ws.Range("M132").Value = -19677.9995
$ws.Range("H136").Value = 3895.4644
$ws.Range("I136").Value = 1149.25
$ws.Range("J136").Value = 7557.0835
$ws.Range("K136").Value = 3447.75
$ws.Range("L136").Value = 22671.2505
$ws.Range("M136").Value = -897.75
$ws.Range("N136").Value = -27771.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1068139.8
$ws.Range("I126").Value = 3533.3333
$ws.Range("J126").Value = 2665049.2
$ws.Range("K126").Value = 10599.9999
$ws.Range("L126").Value = 7995147.600000001
$ws.Range("M126").Value = -8129.999899999999
$ws.Range("N126").Value = -8000087.600000001
$ws.Range("H132").Value = 7579983.5
$ws.Range("I132").Value = 5227.125
$ws.Range("J132").Value = 16669691
$ws.Range("K132").Value = 15681.375
$ws.Range("L132").Value = 50009073
$ws.Range("M132").Value = -13151.375
$ws.Range("N132").Value = -50014133
$ws.Range("H136").Value = 9439.809999999999
$ws.Range("I136").Value = 9063.923000000001
$ws.Range("J136").Value = 10050.625
$ws.Range("K136").Value = 27191.769
$ws.Range("L136").Value = 30151.875
$ws.Range("M136").Value = -24641.769
$ws.Range("N136").Value = -35251.875
